$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(112).RowHeight = 17
$ws.Rows.Item(112).EntireRow.AutoFit()
